# Insert a new data row for the "2026/02/23" date group (16:00 entry)
# before the existing "2026/12/29" block, shifting all subsequent rows
# down by one (dimension grows from A1:D882 to A1:D883).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 841; existing rows 841..882 shift down to 842..883.
$ws.Rows.Item(841).Insert()

# Populate the newly inserted row with the new record.
# Force Text format on the date/weekday columns so Excel doesn't
# auto-convert the date-looking string into a date serial number.
$ws.Cells.Item(841, 1).NumberFormat = "@"
$ws.Cells.Item(841, 1).Value = "2026/02/23"
$ws.Cells.Item(841, 2).Value = "月"
$ws.Cells.Item(841, 3).Value = 16
$ws.Cells.Item(841, 4).Value = 28
